# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-62.
# Rows not listed here already contain 0 and are left unchanged.
$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 3
    25 = 1
    26 = 2
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 3
    35 = 1
    37 = 0
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 2
    53 = 2
    54 = 1
    55 = 0
    56 = 1
    57 = 1
    59 = 1
    60 = 3
    61 = 3
    62 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
